# LMS-2340 Updating basynthec stuff based on results of the leiden meeting.
$wb = $excel.ActiveWorkbook

# --- openbis-metadata sheet: remove the "Strain" property row ---
$wsMeta = $wb.Worksheets.Item("openbis-metadata")
$wsMeta.Rows(3).Delete()

# --- openbis-data sheet: rename header + strain, add rows for extra strains ---
$wsData = $wb.Worksheets.Item("openbis-data")

# Header: "Abs" -> "Strain"
$wsData.Range("A1").Value = "Strain"

# First data row: "OD600" -> "MGP1"
$wsData.Range("A2").Value = "MGP1"

# Copy row 2 (B:U) down into rows 3, 4 and 5 for the additional strains
$wsData.Range("B2:U2").Copy()
$wsData.Range("B3:U3").PasteSpecial()
$wsData.Range("B2:U2").Copy()
$wsData.Range("B4:U4").PasteSpecial()
$wsData.Range("B2:U2").Copy()
$wsData.Range("B5:U5").PasteSpecial()

$wsData.Range("A3").Value = "MGP100"
$wsData.Range("A4").Value = "MGP20"
$wsData.Range("A5").Value = "MGP999"

$wsData.Activate()
$wsData.Range("A13").Select()

# Leave the metadata sheet as the active/selected sheet & row-3 selected
$wsMeta.Activate()
$wsMeta.Range("A3:XFD3").Select()
